# Update column G (header "K") values for rows 2-10 on the active worksheet.
# These values represent a recalculated "K" (strikeouts -> K) statistic that
# replaces the old "Strike#" derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 0
    4  = 4
    5  = 4
    6  = 6
    7  = 2
    8  = 2
    9  = 2
    10 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
